# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap country labels that changed rank order (Casos totales sort) ---
# Canada / Belgica swap (rows 30-31)
$ws.Cells.Item(30, 1).Value = "Belgica"
$ws.Cells.Item(31, 1).Value = "Canada"

# China / Venezuela swap (rows 55-56)
$ws.Cells.Item(55, 1).Value = "Venezuela"
$ws.Cells.Item(56, 1).Value = "China"

# Islas Malvinas / Montserrat swap (rows 216-217)
$ws.Cells.Item(216, 1).Value = "Montserrat"
$ws.Cells.Item(217, 1).Value = "Islas Malvinas"

# --- Update the "last updated" timestamp title ---
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 17 de Octubre de 2020 a las 05:56"

# --- Update numeric data ---

# Row 5 : India
$ws.Cells.Item(5, 2).Value = 7432680
$ws.Cells.Item(5, 3).Value = 2045
$ws.Cells.Item(5, 4).Value = 6524595
$ws.Cells.Item(5, 5).Value = 795053

# Row 30 : now Belgica
$ws.Cells.Item(30, 2).Value = 202151
$ws.Cells.Item(30, 3).Value = 10192
$ws.Cells.Item(30, 4).Value = 20867
$ws.Cells.Item(30, 5).Value = 170925
$ws.Cells.Item(30, 7).Value = 32
$ws.Cells.Item(30, 8).Value = 10359

# Row 31 : now Canada
$ws.Cells.Item(31, 2).Value = 194106
$ws.Cells.Item(31, 4).Value = 163644
$ws.Cells.Item(31, 5).Value = 20740
$ws.Cells.Item(31, 8).Value = 9722

# Row 44 : Kazajistan
$ws.Cells.Item(44, 2).Value = 109302
$ws.Cells.Item(44, 3).Value = 100
$ws.Cells.Item(44, 4).Value = 104921
$ws.Cells.Item(44, 5).Value = 2613

# Row 55 : now Venezuela
$ws.Cells.Item(55, 2).Value = 85758
$ws.Cells.Item(55, 3).Value = 0
$ws.Cells.Item(55, 4).Value = 78294
$ws.Cells.Item(55, 5).Value = 6739
$ws.Cells.Item(55, 8).Value = 725

# Row 56 : now China
$ws.Cells.Item(56, 2).Value = 85659
$ws.Cells.Item(56, 3).Value = 13
$ws.Cells.Item(56, 4).Value = 80766
$ws.Cells.Item(56, 5).Value = 259
$ws.Cells.Item(56, 8).Value = 4634

# Row 187 : Butan
$ws.Cells.Item(187, 4).Value = 298
$ws.Cells.Item(187, 5).Value = 18

# Row 216 : now Montserrat
$ws.Cells.Item(216, 4).Value = 12
$ws.Cells.Item(216, 8).Value = 1

# Row 217 : now Islas Malvinas
$ws.Cells.Item(217, 4).Value = 13
$ws.Cells.Item(217, 8).Value = 0
